# Sprint 6 burndown update — 5th, 6th and a bit of 7th day entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart Sprint6")

# Day 6 (column J) entries
$ws.Range("J6").Value = 0.25
$ws.Range("J7").Value = 0.5
$ws.Range("J16").Value = 0.5
$ws.Range("J21").Value = 2.5

# Day 7 (column K) entries
$ws.Range("K12").Value = 0.25
$ws.Range("K13").Value = 0.25
$ws.Range("K16").Value = 0.5
$ws.Range("K17").Value = 0.25
$ws.Range("K18").Value = 0.5
$ws.Range("K19").Value = 0.5
$ws.Range("K21").Value = 0.5

# Move the active selection, as recorded in the saved view state.
$ws.Range("H25").Select()
